$d = $word.ActiveDocument

# Locate the run of text that needs to be restructured/replaced.
$target = "Weekly meetings on Saturday at 3pm, with additional meetings scheduled as necessary and Daily Scrum meeting everyday at 10pm."
$rng = $d.Content
$found = $rng.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target meeting-schedule sentence"
}

# Re-materialize a clean Range object over the exact same span so that
# InsertXML performs a true replace of that span (rather than an insert
# next to a Find-collapsed range).
$replaceRange = $d.Range($rng.Start, $rng.End)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
'<w:p><w:r><w:t>Daily meetings</w:t></w:r>' +
'<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
'<w:r><w:t xml:space="preserve">every day </w:t></w:r>' +
'<w:r><w:t xml:space="preserve">at </w:t></w:r>' +
'<w:r><w:t>8</w:t></w:r>' +
'<w:r><w:t>pm</w:t></w:r>' +
'<w:r><w:t>, with additional meetings scheduled as necessary</w:t></w:r>' +
'<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
'</w:p>' +
'</w:body></w:document>' +
'</pkg:xmlData></pkg:part></pkg:package>'

$replaceRange.InsertXML($xml)
